$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B3: TOWER value changed
$ws.Range("B3").Value = "CONTROL TOWER"

# B4: PROFILE OVERVIEW text rewritten
$ws.Range("B4").Value = "Data enthusiast with over 6 years of experience in transforming data into clear insights. Specializes in **gathering**, **processing**, and **analyzing** information to support strategic goals. Currently enhancing **Python** and **SQL** skills to develop Business Intelligence solutions. Proven track record in reducing document processing time by 50% and improving data reliability through quality initiatives."

# B6: INDUSTRY EXPERIENCE - removed "Financial Services" line
$ws.Range("B6").Value = "Communications, Media & Technology`nProducts"

# B10: Role_1 - removed the "Managed data quality initiatives..." bullet and tweaked last line
$ws.Range("B10").Value = "Data Analyst  `nResponsible for **collecting** and **analyzing** data on public property in the Polish market, covering 800 companies.  `nContributed to **automation projects** using **Python scripts** and **SQL queries** that reduced document processing time by 50%, optimizing workflows and boosting team efficiency.  `nDeveloped and maintained interactive **Power BI reports** and dashboards to provide actionable insights for stakeholders and support data-driven decision making."

# B11: Role_2 - removed the "Maintained accurate records..." bullet and tweaked a line
$ws.Range("B11").Value = "International Forwarder  `nManaged and supervised **transportation logistics** for freight across multiple European countries, ensuring timely delivery and compliance with regulations.  `nCoordinated communication between clients, drivers, and warehouses to optimize route planning and resolve issues promptly.  `nUtilized **data tracking systems** to monitor shipment progress and identify opportunities to improve delivery efficiency."
